$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('K2').Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276206556/TJaAJ2AJpb-camera-video-segment-driver-1743276209056.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSKDWHZSKL%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T110008Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECIaCXVzLXdlc3QtMiJHMEUCIFxAnevpCvO%2FCXOOdxPWtfhgc4ogmgdm3HSau%2BWYKIABAiEA9rVhPicEkjV2%2FxkiVWXGmXTvUDbx7w6Omsd%2B3F%2Bi0QUq5gMIi%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDMT7qtiGyvNsjboh6yq6Ay2Fxm5K8IGcD21OADW2phrdD2cMfAhpGpl7%2BLDifr7bJbEvYchNFssUSxFOivp%2BCKlIOXSE5ysbZOPpSWsD0avVDaS8hgoM22SMkV1SePof%2BW0zJUIscw03JlYaxj59AYU%2FdFfdXrSqTIDyq217yhEUIMLBDFr%2FxDpXfvPU8ytn%2FilWo2nWiBHabwfKPPNWRkLBfeNHefhFpgA%2FML83qVZgPPkfPQtGxiH9GGIxZZmcpb%2B4IzarFqZz6OZmMV%2F5kBJ37CC0JxjSJJ4Y4CEmOXZfJd5IrV71Tft9xq5kLa4QaXaxFEeTTZMJE%2FikNuzIQAUwcDTl569mJRua8BB7o2Nv4K2SsBLwidB0AwmpyLg7AH8DC4Q9Oyy0SaqR33Ljg3vXshB5UovjvW%2BvKeLQGE04nbLoUWzsxeDPzGzlV%2FPwc9%2FDYQkcIOvjmhUq%2FH5WQ20Lhc8B68cRubVGQ8%2BidryUQsMKHynFTGItDVmyTbc9y%2FXUn%2BNGOI3XpOnn4CqnvdR9kZEEYPwDaSPjuaPZ7j%2Fnb29EVzbTk0Bs9AfoJz8NLDQnm9qS%2BwDQDEbvWIrnubt%2F7WzzYB0dIqEwua6kvwY6pQF3ZhBm2pvjltMFjt4lbEI7SMjBDh53lC7ezc7xRmD5r2QVD7S5CYV2yohJsupVycb1oV8PcTkkOmKso%2Fe5XkOE0ubDLOspH9gBB68p34%2BK6S7wJPalWAlernjUOjK5R2DYRsBkhnjcU8c6jSSwMbwITj%2FYOeumsJbtJ3UrM7zZjql0stsePjZmdkjjWpaE8mdtlNWcuP2jX09eEz4%2Bh9SgXiW2BVI%3D&X-Amz-SignedHeaders=host&response-expires=Sun%2C%2030%20Mar%202025%2019%3A00%3A08%20GMT&X-Amz-Signature=05cffc2ff530c4973ba0d9602c8345ff4c9a2a13c99e32b50d6da77157a94fcd'
$ws.Range('K3').Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276164863/OF6jgx1Lcp-camera-video-segment-driver-1743276167363.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSKDWHZSKL%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T110008Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECIaCXVzLXdlc3QtMiJHMEUCIFxAnevpCvO%2FCXOOdxPWtfhgc4ogmgdm3HSau%2BWYKIABAiEA9rVhPicEkjV2%2FxkiVWXGmXTvUDbx7w6Omsd%2B3F%2Bi0QUq5gMIi%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDMT7qtiGyvNsjboh6yq6Ay2Fxm5K8IGcD21OADW2phrdD2cMfAhpGpl7%2BLDifr7bJbEvYchNFssUSxFOivp%2BCKlIOXSE5ysbZOPpSWsD0avVDaS8hgoM22SMkV1SePof%2BW0zJUIscw03JlYaxj59AYU%2FdFfdXrSqTIDyq217yhEUIMLBDFr%2FxDpXfvPU8ytn%2FilWo2nWiBHabwfKPPNWRkLBfeNHefhFpgA%2FML83qVZgPPkfPQtGxiH9GGIxZZmcpb%2B4IzarFqZz6OZmMV%2F5kBJ37CC0JxjSJJ4Y4CEmOXZfJd5IrV71Tft9xq5kLa4QaXaxFEeTTZMJE%2FikNuzIQAUwcDTl569mJRua8BB7o2Nv4K2SsBLwidB0AwmpyLg7AH8DC4Q9Oyy0SaqR33Ljg3vXshB5UovjvW%2BvKeLQGE04nbLoUWzsxeDPzGzlV%2FPwc9%2FDYQkcIOvjmhUq%2FH5WQ20Lhc8B68cRubVGQ8%2BidryUQsMKHynFTGItDVmyTbc9y%2FXUn%2BNGOI3XpOnn4CqnvdR9kZEEYPwDaSPjuaPZ7j%2Fnb29EVzbTk0Bs9AfoJz8NLDQnm9qS%2BwDQDEbvWIrnubt%2F7WzzYB0dIqEwua6kvwY6pQF3ZhBm2pvjltMFjt4lbEI7SMjBDh53lC7ezc7xRmD5r2QVD7S5CYV2yohJsupVycb1oV8PcTkkOmKso%2Fe5XkOE0ubDLOspH9gBB68p34%2BK6S7wJPalWAlernjUOjK5R2DYRsBkhnjcU8c6jSSwMbwITj%2FYOeumsJbtJ3UrM7zZjql0stsePjZmdkjjWpaE8mdtlNWcuP2jX09eEz4%2Bh9SgXiW2BVI%3D&X-Amz-SignedHeaders=host&response-expires=Sun%2C%2030%20Mar%202025%2019%3A00%3A08%20GMT&X-Amz-Signature=403fd6098bb9b7a0ddc24cfa59746dcd5f33b73bd28a99a88169bcf366619265'
$ws.Range('K4').Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276108738/gAxSKLP636-camera-video-segment-driver-1743276111238.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSKDWHZSKL%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T110008Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECIaCXVzLXdlc3QtMiJHMEUCIFxAnevpCvO%2FCXOOdxPWtfhgc4ogmgdm3HSau%2BWYKIABAiEA9rVhPicEkjV2%2FxkiVWXGmXTvUDbx7w6Omsd%2B3F%2Bi0QUq5gMIi%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDMT7qtiGyvNsjboh6yq6Ay2Fxm5K8IGcD21OADW2phrdD2cMfAhpGpl7%2BLDifr7bJbEvYchNFssUSxFOivp%2BCKlIOXSE5ysbZOPpSWsD0avVDaS8hgoM22SMkV1SePof%2BW0zJUIscw03JlYaxj59AYU%2FdFfdXrSqTIDyq217yhEUIMLBDFr%2FxDpXfvPU8ytn%2FilWo2nWiBHabwfKPPNWRkLBfeNHefhFpgA%2FML83qVZgPPkfPQtGxiH9GGIxZZmcpb%2B4IzarFqZz6OZmMV%2F5kBJ37CC0JxjSJJ4Y4CEmOXZfJd5IrV71Tft9xq5kLa4QaXaxFEeTTZMJE%2FikNuzIQAUwcDTl569mJRua8BB7o2Nv4K2SsBLwidB0AwmpyLg7AH8DC4Q9Oyy0SaqR33Ljg3vXshB5UovjvW%2BvKeLQGE04nbLoUWzsxeDPzGzlV%2FPwc9%2FDYQkcIOvjmhUq%2FH5WQ20Lhc8B68cRubVGQ8%2BidryUQsMKHynFTGItDVmyTbc9y%2FXUn%2BNGOI3XpOnn4CqnvdR9kZEEYPwDaSPjuaPZ7j%2Fnb29EVzbTk0Bs9AfoJz8NLDQnm9qS%2BwDQDEbvWIrnubt%2F7WzzYB0dIqEwua6kvwY6pQF3ZhBm2pvjltMFjt4lbEI7SMjBDh53lC7ezc7xRmD5r2QVD7S5CYV2yohJsupVycb1oV8PcTkkOmKso%2Fe5XkOE0ubDLOspH9gBB68p34%2BK6S7wJPalWAlernjUOjK5R2DYRsBkhnjcU8c6jSSwMbwITj%2FYOeumsJbtJ3UrM7zZjql0stsePjZmdkjjWpaE8mdtlNWcuP2jX09eEz4%2Bh9SgXiW2BVI%3D&X-Amz-SignedHeaders=host&response-expires=Sun%2C%2030%20Mar%202025%2019%3A00%3A08%20GMT&X-Amz-Signature=a6d0823c367af939b2d3accb720d523b7c7e152934b029dd163c47047b10f473'
$ws.Range('K5').Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205341/1743265652740/hq4516Mzws-camera-video-segment-driver-1743265655240.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSKDWHZSKL%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T110008Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECIaCXVzLXdlc3QtMiJHMEUCIFxAnevpCvO%2FCXOOdxPWtfhgc4ogmgdm3HSau%2BWYKIABAiEA9rVhPicEkjV2%2FxkiVWXGmXTvUDbx7w6Omsd%2B3F%2Bi0QUq5gMIi%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDMT7qtiGyvNsjboh6yq6Ay2Fxm5K8IGcD21OADW2phrdD2cMfAhpGpl7%2BLDifr7bJbEvYchNFssUSxFOivp%2BCKlIOXSE5ysbZOPpSWsD0avVDaS8hgoM22SMkV1SePof%2BW0zJUIscw03JlYaxj59AYU%2FdFfdXrSqTIDyq217yhEUIMLBDFr%2FxDpXfvPU8ytn%2FilWo2nWiBHabwfKPPNWRkLBfeNHefhFpgA%2FML83qVZgPPkfPQtGxiH9GGIxZZmcpb%2B4IzarFqZz6OZmMV%2F5kBJ37CC0JxjSJJ4Y4CEmOXZfJd5IrV71Tft9xq5kLa4QaXaxFEeTTZMJE%2FikNuzIQAUwcDTl569mJRua8BB7o2Nv4K2SsBLwidB0AwmpyLg7AH8DC4Q9Oyy0SaqR33Ljg3vXshB5UovjvW%2BvKeLQGE04nbLoUWzsxeDPzGzlV%2FPwc9%2FDYQkcIOvjmhUq%2FH5WQ20Lhc8B68cRubVGQ8%2BidryUQsMKHynFTGItDVmyTbc9y%2FXUn%2BNGOI3XpOnn4CqnvdR9kZEEYPwDaSPjuaPZ7j%2Fnb29EVzbTk0Bs9AfoJz8NLDQnm9qS%2BwDQDEbvWIrnubt%2F7WzzYB0dIqEwua6kvwY6pQF3ZhBm2pvjltMFjt4lbEI7SMjBDh53lC7ezc7xRmD5r2QVD7S5CYV2yohJsupVycb1oV8PcTkkOmKso%2Fe5XkOE0ubDLOspH9gBB68p34%2BK6S7wJPalWAlernjUOjK5R2DYRsBkhnjcU8c6jSSwMbwITj%2FYOeumsJbtJ3UrM7zZjql0stsePjZmdkjjWpaE8mdtlNWcuP2jX09eEz4%2Bh9SgXiW2BVI%3D&X-Amz-SignedHeaders=host&response-expires=Sun%2C%2030%20Mar%202025%2019%3A00%3A08%20GMT&X-Amz-Signature=dd1be111dc93b4f4880e3b9f1a98c595eb776972c49c319de3e5949961919651'
$ws.Range('K6').Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743255999969/Yjq7ZuVo4p-camera-video-segment-driver-1743256002469.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSKDWHZSKL%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T110008Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECIaCXVzLXdlc3QtMiJHMEUCIFxAnevpCvO%2FCXOOdxPWtfhgc4ogmgdm3HSau%2BWYKIABAiEA9rVhPicEkjV2%2FxkiVWXGmXTvUDbx7w6Omsd%2B3F%2Bi0QUq5gMIi%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDMT7qtiGyvNsjboh6yq6Ay2Fxm5K8IGcD21OADW2phrdD2cMfAhpGpl7%2BLDifr7bJbEvYchNFssUSxFOivp%2BCKlIOXSE5ysbZOPpSWsD0avVDaS8hgoM22SMkV1SePof%2BW0zJUIscw03JlYaxj59AYU%2FdFfdXrSqTIDyq217yhEUIMLBDFr%2FxDpXfvPU8ytn%2FilWo2nWiBHabwfKPPNWRkLBfeNHefhFpgA%2FML83qVZgPPkfPQtGxiH9GGIxZZmcpb%2B4IzarFqZz6OZmMV%2F5kBJ37CC0JxjSJJ4Y4CEmOXZfJd5IrV71Tft9xq5kLa4QaXaxFEeTTZMJE%2FikNuzIQAUwcDTl569mJRua8BB7o2Nv4K2SsBLwidB0AwmpyLg7AH8DC4Q9Oyy0SaqR33Ljg3vXshB5UovjvW%2BvKeLQGE04nbLoUWzsxeDPzGzlV%2FPwc9%2FDYQkcIOvjmhUq%2FH5WQ20Lhc8B68cRubVGQ8%2BidryUQsMKHynFTGItDVmyTbc9y%2FXUn%2BNGOI3XpOnn4CqnvdR9kZEEYPwDaSPjuaPZ7j%2Fnb29EVzbTk0Bs9AfoJz8NLDQnm9qS%2BwDQDEbvWIrnubt%2F7WzzYB0dIqEwua6kvwY6pQF3ZhBm2pvjltMFjt4lbEI7SMjBDh53lC7ezc7xRmD5r2QVD7S5CYV2yohJsupVycb1oV8PcTkkOmKso%2Fe5XkOE0ubDLOspH9gBB68p34%2BK6S7wJPalWAlernjUOjK5R2DYRsBkhnjcU8c6jSSwMbwITj%2FYOeumsJbtJ3UrM7zZjql0stsePjZmdkjjWpaE8mdtlNWcuP2jX09eEz4%2Bh9SgXiW2BVI%3D&X-Amz-SignedHeaders=host&response-expires=Sun%2C%2030%20Mar%202025%2019%3A00%3A08%20GMT&X-Amz-Signature=b0a77848a9a229c394ceb32f3cd4d20ec9894e71716f8a7e4b9cc8e94ac3e036'
